# Update the comment under the ReportExplorer folder (rows 17-25).
# For each row in that section, fill in the "Reviewer" (column C) and
# "Status" (column D) cells. Column C needs to pick up the same
# (bold/header-ish) font formatting already used in column B of the same
# block, so we copy formats from the adjacent B cell before setting the
# value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

function Set-ReviewerCell($rowNum) {
    $bCell = $ws.Range("B$rowNum")
    $cCell = $ws.Range("C$rowNum")
    $bCell.Copy()
    $cCell.PasteSpecial($xlPasteFormats)
}

# Row 17 - CreateDashboard.js
Set-ReviewerCell 17
$ws.Range("C17").Value = "Baotong"
$ws.Range("D17").Value = "Done"

# Row 18 - ParameterModel.js (different reviewer, no status note)
Set-ReviewerCell 18
$ws.Range("C18").Value = "Jon"

# Row 19 - ReportExplorer.js
Set-ReviewerCell 19
$ws.Range("C19").Value = "Baotong"
$ws.Range("D19").Value = "Done"

# Row 20 - ReportExplorerContextMenu.js
Set-ReviewerCell 20
$ws.Range("C20").Value = "Baotong"
$ws.Range("D20").Value = "Done"

# Row 21 - ReportExplorerEZ.js
Set-ReviewerCell 21
$ws.Range("C21").Value = "Baotong"
$ws.Range("D21").Value = "Done. transitionToReportViewer(params, urlOption)"

# Row 22 - ReportExplorerSearchFolder.js
Set-ReviewerCell 22
$ws.Range("C22").Value = "Baotong"
$ws.Range("D22").Value = "Done"

# Row 23 - ReportExplorerToolbar.js
Set-ReviewerCell 23
$ws.Range("C23").Value = "Baotong"
$ws.Range("D23").Value = "Done"

# Row 24 - ReportExplorerToolpane.js
Set-ReviewerCell 24
$ws.Range("C24").Value = "Baotong"
$ws.Range("D24").Value = "Done"

# Row 25 - UserSetting.js
Set-ReviewerCell 25
$ws.Range("C25").Value = "Baotong"
$ws.Range("D25").Value = "Done"

# Update the view state to match where the author left the selection.
$ws.Activate()
$ws.Range("C18").Select()
